$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'1112"
$ws.Range("G2").Value = "'9812892189"
$ws.Range("I2").Value = "'2018-11-23"
$ws.Range("O2").Value = "'2018-01-01"
